# Add git to distributions
# Appends a new data row (for the "git" repository) to both worksheets,
# mirroring the existing rows, and updates the active sheet/selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Sheet1 ("all branch segments") - add row 8 for "git"
$ws1.Range("A8").Value = "git"
$ws1.Range("B8").Value = 3.3768472906000002
$ws1.Range("C8").Value = 1
$ws1.Range("D8").Value = 11.374127195
$ws1.Range("E8").Value = 4.5110301992000004
$ws1.Range("F8").Value = 2
$ws1.Range("G8").Value = 9.7064627256999998
$ws1.Range("H8").Value = 1.7696330951000001
$ws1.Range("I8").Value = 1
$ws1.Range("J8").Value = 4.0132459881000004
$ws1.Range("K8").Value = 4.5238809166999996
$ws1.Range("L8").Value = 2
$ws1.Range("M8").Value = 9.7520353179000008
$ws1.Range("N8").Value = 1.7713881890000001
$ws1.Range("O8").Value = 1
$ws1.Range("P8").Value = 4.0132822012
$ws1.Range("Q8").Value = 134.70625401999999
$ws1.Range("R8").Value = 26
$ws1.Range("S8").Value = 705.3084791
$ws1.Range("T8").Value = 51.795459413000003
$ws1.Range("U8").Value = 6
$ws1.Range("V8").Value = 250.72357905000001
$ws1.Range("W8").Value = 186.50171343
$ws1.Range("X8").Value = 36
$ws1.Range("Y8").Value = 848.93752258999996
$ws1.Range("Z8").Value = 59.941029284000003
$ws1.Range("AA8").Value = 19
$ws1.Range("AB8").Value = 246.77768237000001
$ws1.Range("AC8").Value = 154.80316984000001
$ws1.Range("AD8").Value = 26
$ws1.Range("AE8").Value = 1333.1241574000001
$ws1.Range("AF8").Value = 71.892375240999996
$ws1.Range("AG8").Value = 7
$ws1.Range("AH8").Value = 1162.9505340000001
$ws1.Range("AI8").Value = 226.69554507999999
$ws1.Range("AJ8").Value = 37
$ws1.Range("AK8").Value = 2420.6193664000002
$ws1.Range("AL8").Value = 61.30184843
$ws1.Range("AM8").Value = 19.875
$ws1.Range("AN8").Value = 249.0308378
$ws1.Range("AO8").Value = 28.223495395
$ws1.Range("AP8").Value = 6
$ws1.Range("AQ8").Value = 148.22690814000001
$ws1.Range("AR8").Value = 11.953771774
$ws1.Range("AS8").Value = 3
$ws1.Range("AT8").Value = 80.120498069000007
$ws1.Range("AU8").Value = 35.544549154000002
$ws1.Range("AV8").Value = 6
$ws1.Range("AW8").Value = 270.44755229999998
$ws1.Range("AX8").Value = 12.460389075
$ws1.Range("AY8").Value = 3.5
$ws1.Range("AZ8").Value = 80.514722559999996
$ws1.Range("BA8").Value = 1.5515099592999999
$ws1.Range("BB8").Value = 1
$ws1.Range("BC8").Value = 1.3602365662
$ws1.Range("BD8").Value = 1.0126365389
$ws1.Range("BE8").Value = 1
$ws1.Range("BF8").Value = 0.13434073730000001
$ws1.Range("BG8").Value = 46.619348123999998
$ws1.Range("BH8").Value = 0
$ws1.Range("BI8").Value = 630.91710981999995
$ws1.Range("BJ8").Value = 126.32424484000001
$ws1.Range("BK8").Value = 0
$ws1.Range("BL8").Value = 1131.3163784000001

# Sheet2 ("only branch segs gt 1") - add row 8 for "git"
$ws2.Range("A8").Value = "git"
$ws2.Range("B8").Value = 9338
$ws2.Range("C8").Value = 31533
$ws2.Range("D8").Formula = "=B8-4693"
$ws2.Range("E8").Value = 26840
$ws2.Range("F8").Formula = "=D8/B8"
$ws2.Range("G8").Formula = "=E8/C8"
$ws2.Range("H8").Value = 5.7782562000000004
$ws2.Range("I8").Value = 3
$ws2.Range("J8").Value = 15.767991
$ws2.Range("K8").Value = 6.9844995000000001
$ws2.Range("L8").Value = 4
$ws2.Range("M8").Value = 12.232538
$ws2.Range("N8").Value = 1.4733765000000001
$ws2.Range("O8").Value = 1
$ws2.Range("P8").Value = 2.1433053000000002
$ws2.Range("Q8").Value = 7.0103337000000003
$ws2.Range("R8").Value = 4
$ws2.Range("S8").Value = 12.299989
$ws2.Range("T8").Value = 1.4769048
$ws2.Range("U8").Value = 1
$ws2.Range("V8").Value = 2.1439278000000002
$ws2.Range("W8").Value = 227.00388000000001
$ws2.Range("X8").Value = 62
$ws2.Range("Y8").Value = 972.96515999999997
$ws2.Range("Z8").Value = 84.373305000000002
$ws2.Range("AA8").Value = 18
$ws2.Range("AB8").Value = 332.86081000000001
$ws2.Range("AC8").Value = 311.37718000000001
$ws2.Range("AD8").Value = 89
$ws2.Range("AE8").Value = 1154.7445
$ws2.Range("AF8").Value = 56.947972
$ws2.Range("AG8").Value = 23.333300000000001
$ws2.Range("AH8").Value = 194.67671999999999
$ws2.Range("AI8").Value = 267.40537999999998
$ws2.Range("AJ8").Value = 64
$ws2.Range("AK8").Value = 1873.8741
$ws2.Range("AL8").Value = 124.77481
$ws2.Range("AM8").Value = 21
$ws2.Range("AN8").Value = 1643.2081000000001
$ws2.Range("AO8").Value = 392.18018999999998
$ws2.Range("AP8").Value = 93
$ws2.Range("AQ8").Value = 3411.9780000000001
$ws2.Range("AR8").Value = 59.683672999999999
$ws2.Range("AS8").Value = 24.8
$ws2.Range("AT8").Value = 200.39385999999999
$ws2.Range("AU8").Value = 41.981484999999999
$ws2.Range("AV8").Value = 13
$ws2.Range("AW8").Value = 183.5401
$ws2.Range("AX8").Value = 9.2739118999999999
$ws2.Range("AY8").Value = 3.5
$ws2.Range("AZ8").Value = 52.739424999999997
$ws2.Range("BA8").Value = 56.699247
$ws2.Range("BB8").Value = 15
$ws2.Range("BC8").Value = 368.85532999999998
$ws2.Range("BD8").Value = 10.292382
$ws2.Range("BE8").Value = 4
$ws2.Range("BF8").Value = 53.978757000000002
$ws2.Range("BG8").Value = 2.1087191000000001
$ws2.Range("BH8").Value = 2
$ws2.Range("BI8").Value = 1.7612745999999999
$ws2.Range("BJ8").Value = 1.0254037
$ws2.Range("BK8").Value = 1
$ws2.Range("BL8").Value = 0.18963360000000001
$ws2.Range("BM8").Value = 93.720445999999995
$ws2.Range("BN8").Value = 1.24722
$ws2.Range("BO8").Value = 892.13088000000005
$ws2.Range("BP8").Value = 253.95389
$ws2.Range("BQ8").Value = 16.366099999999999
$ws2.Range("BR8").Value = 1594.0009

# Sheet-view bookkeeping: sheet1 becomes the active/selected tab with a
# selection at E14; sheet2 keeps a (different) lingering selection at F13
# and loses the "tabSelected" flag.
$ws2.Range("F13").Select()
$ws1.Activate()
$ws1.Range("E14").Select()
